# The two paragraphs "<id>p112v_1</id>" and "<id>p113r_1</id>" were each
# split across three runs (an opening "<id>" run in Courier New/brownish
# color, a plain middle run with just the bare id text, and a closing
# "</id>" run matching the opening run's formatting). Collapse each back
# into a single run carrying the opening run's formatting, by doing a
# same-text Find/Replace over the whole tagged string: Word's Find engine
# rewrites the matched range as one run using the formatting of the first
# run in the match, merging the three runs into one without altering the
# visible text.
$d = $word.ActiveDocument

$targets = @("<id>p112v_1</id>", "<id>p113r_1</id>")

foreach ($t in $targets) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($t, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $t, 2)
    if (-not $ok) {
        Write-Host "WARNING: target not found:" $t
    }
}
